# Bulk upload template: switch player-identification columns from names to
# passport codes, widen the header columns to fit the new text, and expand
# the instructions sheet with a dedicated "passport codes only" callout.

$wb = $excel.ActiveWorkbook

# --- Sheet "比赛数据" (Match data) ---
$ws1 = $wb.Worksheets.Item("比赛数据")

# Widen columns A:D so the longer "_Passport / ...护照码" headers fit.
# NOTE: the engine's ColumnWidth setter stores (round((input+5/6)*6))/6, i.e.
# it adds ~5/6 of a character (classic Excel cell padding) before rounding to
# the nearest 1/6 of a character, so we back that padding out of the input
# to land as close as possible on the target stored width of 40.83203125.
$ws1.Range("A1:D1").EntireColumn.ColumnWidth = 40.83203125 - 5/6

# Header row now asks for passport codes instead of bare player names
$ws1.Range("A1").Value = "Team_1_Player_1_Passport / 第一队选手一护照码"
$ws1.Range("B1").Value = "Team_1_Player_2_Passport / 第一队选手二护照码"
$ws1.Range("C1").Value = "Team_2_Player_1_Passport / 第二队选手一护照码"
$ws1.Range("D1").Value = "Team_2_Player_2_Passport / 第二队选手二护照码"

# --- Sheet "使用说明" (Instructions) ---
$ws2 = $wb.Worksheets.Item("使用说明")

# Clarify the existing passport-code guidance lines
$ws2.Range("A4").Value = "• 使用选手护照代码（例如：HVGN0BW0, KGLE38K4）- 不是姓名"
$ws2.Range("A5").Value = "• 单打比赛请将第一队选手二护照码和第二队选手二护照码留空"

# Make room for a new "passport codes only" callout block by pushing the
# existing "验证将检查：" section (old rows 14-18) down by 5 rows
$ws2.Rows("14:18").Insert()

$ws2.Range("A14").Value = "重要提示 - 仅使用护照代码："
$ws2.Range("A15").Value = "• 系统要求护照代码，不是选手姓名"
$ws2.Range("A16").Value = "• 每个护照代码为8个字符（字母和数字）"
$ws2.Range("A17").Value = "• 从选手档案或管理面板查找护照代码"
# A18 is left blank as a spacer row (mirrors rows 2, 9 and 13 above)
